# Applies recalculated crafting-leve profit figures (currentAveragePrice /
# NQ+HQ price & profit columns) across the Sheets workbook, per the scheduled
# market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2674.8572
$ws.Range("I64").Value = 2413.6
$ws.Range("J64").Value = 3328
$ws.Range("K64").Value = 2413.6
$ws.Range("L64").Value = 3328
$ws.Range("M64").Value = -2165.6
$ws.Range("N64").Value = -3824
$ws.Range("H67").Value = 2674.8572
$ws.Range("I67").Value = 2413.6
$ws.Range("J67").Value = 3328
$ws.Range("K67").Value = 2413.6
$ws.Range("L67").Value = 3328
$ws.Range("M67").Value = -1555.6
$ws.Range("N67").Value = -5044
$ws.Range("H74").Value = 4027.0833
$ws.Range("I74").Value = 3597.3684
$ws.Range("J74").Value = 5660
$ws.Range("K74").Value = 3597.3684
$ws.Range("L74").Value = 5660
$ws.Range("M74").Value = -2661.3684
$ws.Range("N74").Value = -7532
$ws.Range("H76").Value = 20835988
$ws.Range("J76").Value = 2691.7693
$ws.Range("L76").Value = 2691.7693
$ws.Range("N76").Value = -3321.7693
$ws.Range("H77").Value = 4027.0833
$ws.Range("I77").Value = 3597.3684
$ws.Range("J77").Value = 5660
$ws.Range("K77").Value = 17986.842
$ws.Range("L77").Value = 28300
$ws.Range("M77").Value = -13306.842
$ws.Range("N77").Value = -37660
$ws.Range("H79").Value = 20835988
$ws.Range("J79").Value = 2691.7693
$ws.Range("L79").Value = 2691.7693
$ws.Range("N79").Value = -4875.7693
$ws.Range("H129").Value = 245483.9
$ws.Range("I129").Value = 365.43478
$ws.Range("J129").Value = 558690.8
$ws.Range("K129").Value = 1096.30434
$ws.Range("L129").Value = 1676072.4
$ws.Range("M129").Value = 3903.69566
$ws.Range("N129").Value = -1686072.4
$ws.Range("H132").Value = 3006079
$ws.Range("I132").Value = 3367.0605
$ws.Range("J132").Value = 27778452
$ws.Range("K132").Value = 10101.1815
$ws.Range("L132").Value = 83335356
$ws.Range("M132").Value = -7571.181500000001
$ws.Range("N132").Value = -83340416
$ws.Range("H141").Value = 1228.7106
$ws.Range("I141").Value = 1085.8611
$ws.Range("J141").Value = 3800
$ws.Range("K141").Value = 3257.5833
$ws.Range("L141").Value = 11400
$ws.Range("M141").Value = 1922.4167
$ws.Range("N141").Value = -21760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 120
$ws.Range("I5").Value = 66.666664
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 66.666664
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = 45.333336
$ws.Range("N5").Value = -424
$ws.Range("H32").Value = 7472236.5
$ws.Range("I32").Value = 10275.132
$ws.Range("J32").Value = 35721092
$ws.Range("K32").Value = 10275.132
$ws.Range("L32").Value = 35721092
$ws.Range("M32").Value = -9988.132
$ws.Range("N32").Value = -35721666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 66.666664
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 66.666664
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 48.333336
$ws.Range("N4").Value = -430
$ws.Range("H22").Value = 505.4091
$ws.Range("I22").Value = 505.4091
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 505.4091
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -332.4091
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 1882.9231
$ws.Range("I99").Value = 1425
$ws.Range("J99").Value = 2086.4443
$ws.Range("K99").Value = 1425
$ws.Range("L99").Value = 2086.4443
$ws.Range("M99").Value = 73
$ws.Range("N99").Value = -5082.4443
$ws.Range("H105").Value = 1636.5172
$ws.Range("I105").Value = 1557.1522
$ws.Range("J105").Value = 1940.75
$ws.Range("K105").Value = 1557.1522
$ws.Range("L105").Value = 1940.75
$ws.Range("M105").Value = 189.8478
$ws.Range("N105").Value = -5434.75
$ws.Range("H134").Value = 10505256
$ws.Range("I134").Value = 14286564
$ws.Range("J134").Value = 2233643.5
$ws.Range("K134").Value = 42859692
$ws.Range("L134").Value = 6700930.5
$ws.Range("M134").Value = -42857157
$ws.Range("N134").Value = -6706000.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 13300
$ws.Range("J43").Value = 13300
$ws.Range("L43").Value = 13300
$ws.Range("N43").Value = -13668
$ws.Range("H58").Value = 1111958.4
$ws.Range("I58").Value = 4989.875
$ws.Range("J58").Value = 2674737.5
$ws.Range("K58").Value = 4989.875
$ws.Range("L58").Value = 2674737.5
$ws.Range("M58").Value = -4786.875
$ws.Range("N58").Value = -2675143.5
$ws.Range("H62").Value = 2585.1853
$ws.Range("I62").Value = 2288.125
$ws.Range("K62").Value = 2288.125
$ws.Range("M62").Value = -1664.125
$ws.Range("H65").Value = 2585.1853
$ws.Range("I65").Value = 2288.125
$ws.Range("K65").Value = 11440.625
$ws.Range("M65").Value = -8320.625
$ws.Range("H101").Value = 13300
$ws.Range("J101").Value = 13300
$ws.Range("L101").Value = 13300
$ws.Range("N101").Value = -19790
$ws.Range("H132").Value = 2200.7727
$ws.Range("I132").Value = 1744.6111
$ws.Range("J132").Value = 4253.5
$ws.Range("K132").Value = 5233.8333
$ws.Range("L132").Value = 12760.5
$ws.Range("M132").Value = -2703.8333
$ws.Range("N132").Value = -17820.5
$ws.Range("H134").Value = 980180.1
$ws.Range("I134").Value = 4496.625
$ws.Range("J134").Value = 4449277
$ws.Range("K134").Value = 13489.875
$ws.Range("L134").Value = 13347831
$ws.Range("M134").Value = -10954.875
$ws.Range("N134").Value = -13352901
$ws.Range("H136").Value = 1111958.4
$ws.Range("I136").Value = 4989.875
$ws.Range("J136").Value = 2674737.5
$ws.Range("K136").Value = 14969.625
$ws.Range("L136").Value = 8024212.5
$ws.Range("M136").Value = -12419.625
$ws.Range("N136").Value = -8029312.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7813317
$ws.Range("J131").Value = 978.8298
$ws.Range("L131").Value = 2936.4894
$ws.Range("N131").Value = -13016.4894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 784.9286
$ws.Range("I46").Value = 829.1111
$ws.Range("J46").Value = 705.4
$ws.Range("K46").Value = 829.1111
$ws.Range("L46").Value = 705.4
$ws.Range("M46").Value = -641.1111
$ws.Range("N46").Value = -1081.4
$ws.Range("H82").Value = 4075.85
$ws.Range("I82").Value = 899.7
$ws.Range("J82").Value = 7252
$ws.Range("K82").Value = 899.7
$ws.Range("L82").Value = 7252
$ws.Range("M82").Value = -538.7
$ws.Range("N82").Value = -7974
$ws.Range("H85").Value = 4075.85
$ws.Range("I85").Value = 899.7
$ws.Range("J85").Value = 7252
$ws.Range("K85").Value = 899.7
$ws.Range("L85").Value = 7252
$ws.Range("M85").Value = 348.3
$ws.Range("N85").Value = -9748
